# Generate Report for Handback
# Adds the handback row for a1a03878-39ca-405f-8986-3a5e061c1919.md
# to the "Overview", "zh-cn" and "de-de" sheets/tables.

$wb = $excel.ActiveWorkbook

$guid = "a1a03878-39ca-405f-8986-3a5e061c1919"
$mdName = "$guid.md"
$mdDisplay = "e2e\$guid.md"

# ---------------------------------------------------------------
# Sheet "Overview" (table3 / displayName "Overview")
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()
$rngOverview = $rowOverview.Range

$rngOverview.Cells.Item(1, 1).Value = $mdName
$wsOverview.Hyperlinks.Add(
    $rngOverview.Cells.Item(1, 2),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/52a9fc47488a062bac71fdea8eb806c0c6d9ee93/$mdDisplay",
    "",
    "",
    $mdDisplay
)
$rngOverview.Cells.Item(1, 3).Value = ".md"
$rngOverview.Cells.Item(1, 5).Value = "Handed back: in sync with en-US"
$rngOverview.Cells.Item(1, 6).Value = "Handed back: in sync with en-US"
$rngOverview.Cells.Item(1, 7).Value = "2016-08-19 06:43:33"
$rngOverview.Cells.Item(1, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------
# Sheet "zh-cn" (table1)
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$rowZhCn = $loZhCn.ListRows.Add()
$rngZhCn = $rowZhCn.Range

$xlfZhCn = "$guid.8587886e6ddbbe9ca1ea4e3caf98ea221ee3ea18.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add(
    $rngZhCn.Cells.Item(1, 1),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/52a9fc47488a062bac71fdea8eb806c0c6d9ee93/$mdDisplay",
    "",
    "",
    $mdName
)
$rngZhCn.Cells.Item(1, 2).Value = ".md"
$rngZhCn.Cells.Item(1, 3).Value = "Handed back: in sync with en-US"
$rngZhCn.Cells.Item(1, 4).Value = "e2e"
$rngZhCn.Cells.Item(1, 5).Value = "ht"
$rngZhCn.Cells.Item(1, 6).Value = "True"
$rngZhCn.Cells.Item(1, 7).Value = $xlfZhCn
$rngZhCn.Cells.Item(1, 8).Value = "2016-08-19 06:43:28"
$rngZhCn.Cells.Item(1, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Hyperlinks.Add(
    $rngZhCn.Cells.Item(1, 9),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9bb37d3f6c9805aa190d44918724911b5216ba10/$mdDisplay",
    "",
    "",
    $mdName
)
$rngZhCn.Cells.Item(1, 10).Value = $xlfZhCn
$rngZhCn.Cells.Item(1, 11).Value = "2016-08-19 06:43:45"
$rngZhCn.Cells.Item(1, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$rngZhCn.Cells.Item(1, 12).Value = ""
$rngZhCn.Cells.Item(1, 13).Value = "True"
$rngZhCn.Cells.Item(1, 14).Value = ""
$rngZhCn.Cells.Item(1, 15).Value = "False"
$rngZhCn.Cells.Item(1, 16).Value = ""

# ---------------------------------------------------------------
# Sheet "de-de" (table2)
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$rowDeDe = $loDeDe.ListRows.Add()
$rngDeDe = $rowDeDe.Range

$xlfDeDe = "$guid.8587886e6ddbbe9ca1ea4e3caf98ea221ee3ea18.de-de.xlf"

$wsDeDe.Hyperlinks.Add(
    $rngDeDe.Cells.Item(1, 1),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/52a9fc47488a062bac71fdea8eb806c0c6d9ee93/$mdDisplay",
    "",
    "",
    $mdName
)
$rngDeDe.Cells.Item(1, 2).Value = ".md"
$rngDeDe.Cells.Item(1, 3).Value = "Handed back: in sync with en-US"
$rngDeDe.Cells.Item(1, 4).Value = "e2e"
$rngDeDe.Cells.Item(1, 5).Value = "ht"
$rngDeDe.Cells.Item(1, 6).Value = "True"
$rngDeDe.Cells.Item(1, 7).Value = $xlfDeDe
$rngDeDe.Cells.Item(1, 8).Value = "2016-08-19 06:43:33"
$rngDeDe.Cells.Item(1, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Hyperlinks.Add(
    $rngDeDe.Cells.Item(1, 9),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/fc69a5d3e308a52b88b61f36f3eb4a78538be279/$mdDisplay",
    "",
    "",
    $mdName
)
$rngDeDe.Cells.Item(1, 10).Value = $xlfDeDe
$rngDeDe.Cells.Item(1, 11).Value = "2016-08-19 06:43:52"
$rngDeDe.Cells.Item(1, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$rngDeDe.Cells.Item(1, 12).Value = ""
$rngDeDe.Cells.Item(1, 13).Value = "True"
$rngDeDe.Cells.Item(1, 14).Value = ""
$rngDeDe.Cells.Item(1, 15).Value = "False"
$rngDeDe.Cells.Item(1, 16).Value = ""

Write-Output "Added handback row for $mdName to Overview, zh-cn and de-de sheets."
